$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.932.92"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "2.997.99"
$ws.Range("E3").Value = "  -5.58%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.80"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.95"
$ws.Range("E6").Value = "  -6.30%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "2.990.33"
$ws.Range("E8").Value = "  -5.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("E10").Value = "  -6.36%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("E12").Value = "  -4.94%  "
$ws.Range("E13").Value = "  -6.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.86"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "3.492.99"
$ws.Range("E16").Value = "  -5.51%  "
$ws.Range("D17").Value = "61.000.47"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "3.003.31"
$ws.Range("E18").Value = "  -5.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  -7.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.40"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  -5.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.662"
$ws.Range("E22").Value = "  -7.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  -7.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.83"
$ws.Range("E25").Value = "  -5.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  -7.59%  "
$ws.Range("E29").Value = "  -8.26%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.49"
$ws.Range("E30").Value = "  -6.75%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.87"
$ws.Range("E31").Value = "  -7.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -10.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0937"
$ws.Range("E33").Value = "  -9.63%  "
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.954"
$ws.Range("E35").Value = "  -8.57%  "
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.07"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "0.0₃0665"
$ws.Range("E38").Value = "  -6.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0358"
$ws.Range("E39").Value = "  -8.02%  "
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("E41").Value = "  -4.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "375.67"
$ws.Range("E42").Value = "  -7.37%  "
$ws.Range("D43").Value = "2.683.73"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("E44").Value = "  -9.11%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.93"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.234"
$ws.Range("E47").Value = "  -7.65%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.46"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  -8.23%  "
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.25"
$ws.Range("E51").Value = "  -9.96%  "
